# Fix element order inside <w:rPr> of several "Tok" character styles so
# that <w:b/> / <w:i/> precede <w:color/>, matching the wml.xsd schema
# sequence (CT_RPr lists rStyle, rFonts, b, bCs, i, iCs, ... color, ...).
#
# Re-assigning the Font.Bold / Font.Italic properties on each Style
# causes the document engine to re-emit the run-properties in the
# canonical schema order, which moves <w:b/> / <w:i/> ahead of
# <w:color/> without altering any values.

$d = $word.ActiveDocument

# Styles whose <w:rPr> only contains <w:color/> and <w:b/>
$boldOnly = @(
    "KeywordTok",
    "ImportTok",
    "ControlFlowTok",
    "AlertTok",
    "ErrorTok"
)

# Styles whose <w:rPr> only contains <w:color/> and <w:i/>
$italicOnly = @(
    "CommentTok",
    "DocumentationTok"
)

# Styles whose <w:rPr> contains <w:color/>, <w:b/> and <w:i/>
$boldItalic = @(
    "AnnotationTok",
    "CommentVarTok",
    "InformationTok",
    "WarningTok"
)

foreach ($styleName in $boldOnly) {
    $s = $d.Styles($styleName)
    $s.Font.Bold = $s.Font.Bold
}

foreach ($styleName in $italicOnly) {
    $s = $d.Styles($styleName)
    $s.Font.Italic = $s.Font.Italic
}

foreach ($styleName in $boldItalic) {
    $s = $d.Styles($styleName)
    $s.Font.Bold = $s.Font.Bold
    $s.Font.Italic = $s.Font.Italic
}
